$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '66.963.95'
$c.Style = 'Normal'
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  +6.10%  '
$c.Style = 'Normal'
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '3.855.06'
$c.Style = 'Normal'
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  +7.88%  '
$c.Style = 'Normal'
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  -0.08%  '
$c.Style = 'Normal'
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '424.27'
$c.Style = 'Normal'
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  +6.02%  '
$c.Style = 'Normal'
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '130.21'
$c.Style = 'Normal'
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  -0.57%  '
$c.Style = 'Normal'
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '3.848.77'
$c.Style = 'Normal'
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  +7.66%  '
$c.Style = 'Normal'
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.610'
$c.Style = 'Normal'
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  +0.39%  '
$c.Style = 'Normal'
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  -0.13%  '
$c.Style = 'Normal'
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.725'
$c.Style = 'Normal'
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  +2.37%  '
$c.Style = 'Normal'
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.159'
$c.Style = 'Normal'
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +1.78%  '
$c.Style = 'Normal'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.0000339'
$c.Style = 'Normal'
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  +9.93%  '
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '40.83'
$c.Style = 'Normal'
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  +0.17%  '
$c.Style = 'Normal'
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  +6.33%  '
$c.Style = 'Normal'
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '4.468.18'
$c.Style = 'Normal'
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  +7.98%  '
$c.Style = 'Normal'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '15.83'
$c.Style = 'Normal'
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  +20.95%  '
$c.Style = 'Normal'
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '3.844.48'
$c.Style = 'Normal'
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  +7.57%  '
$c.Style = 'Normal'
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  -0.52%  '
$c.Style = 'Normal'
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '19.96'
$c.Style = 'Normal'
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  +1.97%  '
$c.Style = 'Normal'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '67.205.91'
$c.Style = 'Normal'
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  +6.39%  '
$c.Style = 'Normal'
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  +2.01%  '
$c.Style = 'Normal'
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '414.72'
$c.Style = 'Normal'
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  +0.61%  '
$c.Style = 'Normal'
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '14.95'
$c.Style = 'Normal'
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  +1.22%  '
$c.Style = 'Normal'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '84.33'
$c.Style = 'Normal'
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  +0.46%  '
$c.Style = 'Normal'
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  +3.62%  '
$c.Style = 'Normal'
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '37.65'
$c.Style = 'Normal'
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  +7.81%  '
$c.Style = 'Normal'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '9.98'
$c.Style = 'Normal'
$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  +8.48%  '
$c.Style = 'Normal'
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '3.25'
$c.Style = 'Normal'
$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  +3.91%  '
$c.Style = 'Normal'
$c = $ws.Range('B29')
$c.NumberFormat = '@'
$c.Value = 'LEO'
$c.Style = 'Normal'
$c = $ws.Range('C29')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$c.Style = 'Normal'
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '5.30'
$c.Style = 'Normal'
$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  +4.17%  '
$c.Style = 'Normal'
$c = $ws.Range('B30')
$c.NumberFormat = '@'
$c.Value = 'RenderToken'
$c.Style = 'Normal'
$c = $ws.Range('C30')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c.Style = 'Normal'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '9.21'
$c.Style = 'Normal'
$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  +35.42%  '
$c.Style = 'Normal'
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '734.25'
$c.Style = 'Normal'
$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  +10.08%  '
$c.Style = 'Normal'
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '13.19'
$c.Style = 'Normal'
$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  +6.02%  '
$c.Style = 'Normal'
$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  +7.99%  '
$c.Style = 'Normal'
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '2.70'
$c.Style = 'Normal'
$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  +0.76%  '
$c.Style = 'Normal'
$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  -0.07%  '
$c.Style = 'Normal'
$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  -2.98%  '
$c.Style = 'Normal'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '38.73'
$c.Style = 'Normal'
$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  -2.91%  '
$c.Style = 'Normal'
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  +28.19%  '
$c.Style = 'Normal'
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  +0.38%  '
$c.Style = 'Normal'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.0₃0741'
$c.Style = 'Normal'
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  +18.95%  '
$c.Style = 'Normal'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.0464'
$c.Style = 'Normal'
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  +2.08%  '
$c.Style = 'Normal'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '2.90'
$c.Style = 'Normal'
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  +2.78%  '
$c.Style = 'Normal'
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  +0.87%  '
$c.Style = 'Normal'
$c = $ws.Range('B44')
$c.NumberFormat = '@'
$c.Value = 'Stellar'
$c.Style = 'Normal'
$c = $ws.Range('C44')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c.Style = 'Normal'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.135'
$c.Style = 'Normal'
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  -2.07%  '
$c.Style = 'Normal'
$c = $ws.Range('B45')
$c.NumberFormat = '@'
$c.Value = 'LidoDAOToken'
$c.Style = 'Normal'
$c = $ws.Range('C45')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c.Style = 'Normal'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '3.36'
$c.Style = 'Normal'
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  +4.11%  '
$c.Style = 'Normal'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.317'
$c.Style = 'Normal'
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  +11.33%  '
$c.Style = 'Normal'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '3.13'
$c.Style = 'Normal'
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  +0.22%  '
$c.Style = 'Normal'
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  +1.83%  '
$c.Style = 'Normal'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '2.82'
$c.Style = 'Normal'
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  +2.13%  '
$c.Style = 'Normal'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '139.71'
$c.Style = 'Normal'
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  -2.75%  '
$c.Style = 'Normal'
$c = $ws.Range('B51')
$c.NumberFormat = '@'
$c.Value = 'WEMIXToken'
$c.Style = 'Normal'
$c = $ws.Range('C51')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c.Style = 'Normal'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '2.56'
$c.Style = 'Normal'
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  +2.14%  '
$c.Style = 'Normal'
